$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'286.48"
$ws.Range("E2").Value = "'2.58%"
$ws.Range("G2").Value = "'21"
$ws.Range("D3").Value = "'28.77"
$ws.Range("E3").Value = "'4.73%"
$ws.Range("G3").Value = "'21"
$ws.Range("D4").Value = "'5.055"
$ws.Range("E4").Value = "'4.58%"
$ws.Range("G4").Value = "'21"
$ws.Range("D5").Value = "'0.06711"
$ws.Range("E5").Value = "'5.23%"
$ws.Range("G5").Value = "'21"
$ws.Range("D6").Value = "'7.327"
$ws.Range("E6").Value = "'4.39%"
$ws.Range("G6").Value = "'21"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'3.393"
$ws.Range("E7").Value = "'2.02%"
$ws.Range("G7").Value = "'21"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "'1.380"
$ws.Range("E8").Value = "'6.99%"
$ws.Range("G8").Value = "'21"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9414"
$ws.Range("E9").Value = "'5.38%"
$ws.Range("G9").Value = "'21"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1568"
$ws.Range("E10").Value = "'2.13%"
$ws.Range("G10").Value = "'21"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.06812"
$ws.Range("E11").Value = "'12.87%"
$ws.Range("G11").Value = "'21"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07594"
$ws.Range("E12").Value = "'1.15%"
$ws.Range("G12").Value = "'21"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.02951"
$ws.Range("E13").Value = "'0.48%"
$ws.Range("G13").Value = "'21"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09008"
$ws.Range("E14").Value = "'0.15%"
$ws.Range("G14").Value = "'21"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001590"
$ws.Range("E15").Value = "'1.89%"
$ws.Range("G15").Value = "'21"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "'0.04481"
$ws.Range("E16").Value = "'1.84%"
$ws.Range("G16").Value = "'21"
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D17").Value = "'0.0006470"
$ws.Range("E17").Value = "'1.16%"
$ws.Range("G17").Value = "'21"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "'0.006318"
$ws.Range("E18").Value = "'5.24%"
$ws.Range("G18").Value = "'21"
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D19").Value = "'3.449"
$ws.Range("E19").Value = "'-0.98%"
$ws.Range("G19").Value = "'21"
$ws.Range("D20").Value = "'2.253"
$ws.Range("E20").Value = "'1.07%"
$ws.Range("G20").Value = "'21"
$ws.Range("D21").Value = "'0.3210"
$ws.Range("E21").Value = "'2.00%"
$ws.Range("G21").Value = "'21"
$ws.Range("E22").Value = "'-2.96%"
$ws.Range("G22").Value = "'21"
$ws.Range("E23").Value = "'3.91%"
$ws.Range("G23").Value = "'21"
$ws.Range("E24").Value = "'3.11%"
$ws.Range("G24").Value = "'21"
$ws.Range("D25").Value = "'0.001178"
$ws.Range("E25").Value = "'0.21%"
$ws.Range("G25").Value = "'21"
$ws.Range("D26").Value = "'0.004494"
$ws.Range("E26").Value = "'4.94%"
$ws.Range("G26").Value = "'21"
$ws.Range("D27").Value = "'0.0001247"
$ws.Range("E27").Value = "'5.67%"
$ws.Range("G27").Value = "'21"
$ws.Range("D28").Value = "'0.0001614"
$ws.Range("E28").Value = "'-2.31%"
$ws.Range("G28").Value = "'21"
$ws.Range("G29").Value = "'21"
$ws.Range("G30").Value = "'21"
$ws.Range("G31").Value = "'21"
$ws.Range("G32").Value = "'21"
$ws.Range("G33").Value = "'21"
$ws.Range("G34").Value = "'21"
$ws.Range("G35").Value = "'21"
$ws.Range("G36").Value = "'21"
$ws.Range("G37").Value = "'21"
$ws.Range("G38").Value = "'21"
$ws.Range("G39").Value = "'21"
$ws.Range("D40").Value = "'0.04202"
$ws.Range("E40").Value = "'3.16%"
$ws.Range("G40").Value = "'21"
$ws.Range("D41").Value = "'0.006725"
$ws.Range("E41").Value = "'2.17%"
$ws.Range("G41").Value = "'21"
$ws.Range("D42").Value = "'0.1257"
$ws.Range("E42").Value = "'-10.93%"
$ws.Range("G42").Value = "'21"
$ws.Range("D43").Value = "'0.002015"
$ws.Range("E43").Value = "'-3.13%"
$ws.Range("G43").Value = "'21"
$ws.Range("D44").Value = "'0.01229"
$ws.Range("E44").Value = "'11.77%"
$ws.Range("G44").Value = "'21"
$ws.Range("D45").Value = "'0.00005646"
$ws.Range("E45").Value = "'1.82%"
$ws.Range("G45").Value = "'21"
$ws.Range("E46").Value = "'25.93%"
$ws.Range("G46").Value = "'21"
$ws.Range("D47").Value = "'0.01304"
$ws.Range("E47").Value = "'-29.44%"
$ws.Range("G47").Value = "'21"
$ws.Range("G48").Value = "'21"
$ws.Range("G49").Value = "'21"
$ws.Range("G50").Value = "'21"
$ws.Range("G51").Value = "'21"

# Reset styles on the touched range so forced-text cells (apostrophe-prefixed)
# do not retain a quotePrefix style, matching the source formatting (no explicit style).
$ws.Range("A2:G51").Style = "Normal"
